$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 4.646600000000001
$ws.Range("B10").Value = 8.672000000000004
$ws.Range("B12").Value = 6.5923
$ws.Range("D13").Value = -7.548400000000003
$ws.Range("B18").Value = 6.412400000000003
$ws.Range("B25").Value = 5.456999999999997
